# Realestate Update resale numbers 2023-06-02 22:24
# Append a new data row (row 17) to the CityResaleNum sheet with the
# latest resale numbers snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 17

# Columns A-D (Date/Time/Weekday/Week) are stored as text in this sheet,
# so force text formatting before assigning to avoid Excel auto-converting
# the date/time-looking strings into numeric date/time serials.
$textRange = $ws.Range("A" + $row + ":D" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-02"
$ws.Cells.Item($row, 2).Value = "22:21:34"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "22"

# Drop the temporary "Text" number-format override so the new cells end up
# with the same default (unstyled) formatting as the rest of the sheet.
$textRange.ClearFormats()

$ws.Cells.Item($row, 5).Value = 120654
$ws.Cells.Item($row, 6).Value = 133793
$ws.Cells.Item($row, 7).Value = 158797
$ws.Cells.Item($row, 8).Value = 130163
$ws.Cells.Item($row, 9).Value = 174389
$ws.Cells.Item($row, 10).Value = 112238
$ws.Cells.Item($row, 11).Value = 199466
$ws.Cells.Item($row, 12).Value = 218299
$ws.Cells.Item($row, 13).Value = 171424
$ws.Cells.Item($row, 14).Value = 118670
$ws.Cells.Item($row, 15).Value = 38052
$ws.Cells.Item($row, 16).Value = 34846
$ws.Cells.Item($row, 17).Value = 50068
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36899
$ws.Cells.Item($row, 20).Value = -1
